# Weekly NYPD CompStat report refresh: new crime data collected.
# Updates the report period text and all the Week-to-Date / 28-Day / Year-to-Date /
# 2-Year crime-count and percent-change figures for rows 14-31.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Update the rich-text header captions (issue number and report week dates)
# ---------------------------------------------------------------------------
# A8  holds "Volume 32   Number  28"  -> bump the issue number to 29
# C9  holds "Report Covering the Week  7/7/2025  Through  7/13/2025" -> next week
$ws.Range("A8").Value = "Volume 32   Number  29"
$ws.Range("C9").Value = "Report Covering the Week  7/14/2025  Through  7/20/2025"

# ---------------------------------------------------------------------------
# 2) Update the numeric crime-statistics grid (columns C:N, rows 14-31)
# ---------------------------------------------------------------------------
$numericUpdates = @(
  @{ Cell = "N14"; Value = -94.736842105263 },

  @{ Cell = "C15"; Value = 2 },
  @{ Cell = "D15"; Value = 1 },
  @{ Cell = "E15"; Value = 100 },
  @{ Cell = "F15"; Value = 5 },
  @{ Cell = "G15"; Value = 7 },
  @{ Cell = "H15"; Value = -28.571428571428 },
  @{ Cell = "I15"; Value = 24 },
  @{ Cell = "J15"; Value = 23 },
  @{ Cell = "K15"; Value = 4.347826086956 },
  @{ Cell = "L15"; Value = 33.333333333333 },
  @{ Cell = "M15"; Value = 118.181818181818 },
  @{ Cell = "N15"; Value = 100 },

  @{ Cell = "C16"; Value = 8 },
  @{ Cell = "D16"; Value = 12 },
  @{ Cell = "E16"; Value = -33.333333333333 },
  @{ Cell = "F16"; Value = 39 },
  @{ Cell = "G16"; Value = 53 },
  @{ Cell = "H16"; Value = -26.415094339622 },
  @{ Cell = "I16"; Value = 218 },
  @{ Cell = "J16"; Value = 288 },
  @{ Cell = "K16"; Value = -24.305555555555 },
  @{ Cell = "L16"; Value = 1.869158878504 },
  @{ Cell = "M16"; Value = 17.204301075268 },
  @{ Cell = "N16"; Value = -75.227272727272 },

  @{ Cell = "C17"; Value = 13 },
  @{ Cell = "D17"; Value = 22 },
  @{ Cell = "E17"; Value = -40.909090909090 },
  @{ Cell = "F17"; Value = 38 },
  @{ Cell = "G17"; Value = 86 },
  @{ Cell = "H17"; Value = -55.813953488372 },
  @{ Cell = "I17"; Value = 323 },
  @{ Cell = "J17"; Value = 467 },
  @{ Cell = "K17"; Value = -30.835117773019 },
  @{ Cell = "L17"; Value = -18.227848101265 },
  @{ Cell = "M17"; Value = 118.243243243243 },
  @{ Cell = "N17"; Value = 14.539007092198 },

  @{ Cell = "C18"; Value = 1 },
  @{ Cell = "D18"; Value = 6 },
  @{ Cell = "E18"; Value = -83.333333333333 },
  @{ Cell = "F18"; Value = 13 },
  @{ Cell = "G18"; Value = 22 },
  @{ Cell = "H18"; Value = -40.909090909090 },
  @{ Cell = "I18"; Value = 100 },
  @{ Cell = "J18"; Value = 148 },
  @{ Cell = "K18"; Value = -32.432432432432 },
  @{ Cell = "L18"; Value = -12.280701754386 },
  @{ Cell = "M18"; Value = -40.119760479041 },
  @{ Cell = "N18"; Value = -92.125984251968 },

  @{ Cell = "C19"; Value = 10 },
  @{ Cell = "D19"; Value = 28 },
  @{ Cell = "E19"; Value = -64.285714285714 },
  @{ Cell = "F19"; Value = 60 },
  @{ Cell = "G19"; Value = 89 },
  @{ Cell = "H19"; Value = -32.584269662921 },
  @{ Cell = "I19"; Value = 414 },
  @{ Cell = "J19"; Value = 613 },
  @{ Cell = "K19"; Value = -32.463295269168 },
  @{ Cell = "L19"; Value = -30.185497470489 },
  @{ Cell = "M19"; Value = 43.252595155709 },
  @{ Cell = "N19"; Value = -36.209553158705 },

  @{ Cell = "C20"; Value = 2 },
  @{ Cell = "D20"; Value = 7 },
  @{ Cell = "E20"; Value = -71.428571428571 },
  @{ Cell = "F20"; Value = 18 },
  @{ Cell = "G20"; Value = 28 },
  @{ Cell = "H20"; Value = -35.714285714285 },
  @{ Cell = "I20"; Value = 125 },
  @{ Cell = "J20"; Value = 153 },
  @{ Cell = "K20"; Value = -18.300653594771 },
  @{ Cell = "L20"; Value = -29.775280898876 },
  @{ Cell = "M20"; Value = 42.045454545454 },
  @{ Cell = "N20"; Value = -89.397794741306 },

  @{ Cell = "C21"; Value = 36 },
  @{ Cell = "D21"; Value = 76 },
  @{ Cell = "E21"; Value = -52.631578947368 },
  @{ Cell = "F21"; Value = 173 },
  @{ Cell = "G21"; Value = 285 },
  @{ Cell = "H21"; Value = -39.298245614035 },
  @{ Cell = "I21"; Value = 1205 },
  @{ Cell = "J21"; Value = 1694 },
  @{ Cell = "K21"; Value = -28.866587957497 },
  @{ Cell = "L21"; Value = -20.409511228533 },
  @{ Cell = "M21"; Value = 35.393258426966 },
  @{ Cell = "N21"; Value = -71.917967839664 },

  @{ Cell = "C22"; Value = 2 },
  @{ Cell = "F22"; Value = 5 },
  @{ Cell = "G22"; Value = 4 },
  @{ Cell = "H22"; Value = 25 },
  @{ Cell = "I22"; Value = 30 },
  @{ Cell = "J22"; Value = 28 },
  @{ Cell = "K22"; Value = 7.142857142857 },
  @{ Cell = "L22"; Value = 15.384615384615 },
  @{ Cell = "M22"; Value = 66.666666666666 },

  @{ Cell = "C24"; Value = 51 },
  @{ Cell = "D24"; Value = 71 },
  @{ Cell = "E24"; Value = -28.169014084507 },
  @{ Cell = "F24"; Value = 202 },
  @{ Cell = "G24"; Value = 244 },
  @{ Cell = "H24"; Value = -17.213114754098 },
  @{ Cell = "I24"; Value = 1339 },
  @{ Cell = "J24"; Value = 1787 },
  @{ Cell = "K24"; Value = -25.069949636261 },
  @{ Cell = "L24"; Value = -12.768729641693 },
  @{ Cell = "M24"; Value = 46.659364731653 },

  @{ Cell = "C25"; Value = 35 },
  @{ Cell = "D25"; Value = 54 },
  @{ Cell = "E25"; Value = -35.185185185185 },
  @{ Cell = "F25"; Value = 147 },
  @{ Cell = "G25"; Value = 180 },
  @{ Cell = "H25"; Value = -18.333333333333 },
  @{ Cell = "I25"; Value = 989 },
  @{ Cell = "J25"; Value = 1424 },
  @{ Cell = "K25"; Value = -30.547752808988 },
  @{ Cell = "L25"; Value = -12.709620476610 },

  @{ Cell = "C26"; Value = 12 },
  @{ Cell = "D26"; Value = 31 },
  @{ Cell = "E26"; Value = -61.290322580645 },
  @{ Cell = "F26"; Value = 89 },
  @{ Cell = "G26"; Value = 121 },
  @{ Cell = "H26"; Value = -26.446280991735 },
  @{ Cell = "I26"; Value = 562 },
  @{ Cell = "J26"; Value = 739 },
  @{ Cell = "K26"; Value = -23.951285520974 },
  @{ Cell = "L26"; Value = -10.509554140127 },
  @{ Cell = "M26"; Value = 61.031518624641 },

  @{ Cell = "C27"; Value = 2 },
  @{ Cell = "D27"; Value = 2 },
  @{ Cell = "E27"; Value = 0 },
  @{ Cell = "F27"; Value = 7 },
  @{ Cell = "G27"; Value = 9 },
  @{ Cell = "H27"; Value = -22.222222222222 },
  @{ Cell = "I27"; Value = 31 },
  @{ Cell = "J27"; Value = 34 },
  @{ Cell = "K27"; Value = -8.823529411764 },
  @{ Cell = "L27"; Value = 14.814814814814 },

  @{ Cell = "C28"; Value = 3 },
  @{ Cell = "D28"; Value = 1 },
  @{ Cell = "E28"; Value = 200 },
  @{ Cell = "F28"; Value = 9 },
  @{ Cell = "G28"; Value = 14 },
  @{ Cell = "H28"; Value = -35.714285714285 },
  @{ Cell = "I28"; Value = 66 },
  @{ Cell = "J28"; Value = 75 },
  @{ Cell = "K28"; Value = -12 },
  @{ Cell = "L28"; Value = -15.384615384615 },

  @{ Cell = "N29"; Value = -86.842105263157 },

  @{ Cell = "N30"; Value = -91.176470588235 },

  @{ Cell = "L31"; Value = -14.285714285714 }
)

foreach ($u in $numericUpdates) {
  $ws.Range($u.Cell).Value = $u.Value
}

# ---------------------------------------------------------------------------
# 3) Cells that flip from a numeric figure to the "not applicable" placeholder
#    text ("0" / "***.*"), which are stored as shared strings with a
#    dedicated number format (style index 13). Copy that exact format from a
#    cell that already carries it, then apply the matching placeholder text.
# ---------------------------------------------------------------------------
$placeholderTargets = @(
  @{ Cell = "D22"; Text = "'0" },
  @{ Cell = "E22"; Text = "'***.*" },
  @{ Cell = "G31"; Text = "'0" },
  @{ Cell = "H31"; Text = "'***.*" }
)

foreach ($p in $placeholderTargets) {
  # Put the text in first (this may temporarily pick up a generic text style) ...
  $ws.Range($p.Cell).Value = $p.Text
  # ... then copy the canonical placeholder formatting (style 13) from C14
  # (the "0" placeholder) or E14 (the "***.*" placeholder) over it, without
  # touching the value we just set.
  if ($p.Text -eq "'0") {
    $ws.Range("C14").Copy()
  } else {
    $ws.Range("E14").Copy()
  }
  $ws.Range($p.Cell).PasteSpecial(-4122)
}
